# Auto-generated edit script applying the diff to Sheets/Exodus_Profits.xlsx
# Updates specific H/I/J/K/L/M/N numeric cells across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 50664.332
$ws.Range("I57").Value = 40000
$ws.Range("K57").Value = 120000
$ws.Range("M57").Value = -119501
$ws.Range("H70").Value = 1574.75
$ws.Range("J70").Value = 1633.1666
$ws.Range("L70").Value = 4899.4998
$ws.Range("N70").Value = -5439.4998
$ws.Range("H73").Value = 1574.75
$ws.Range("J73").Value = 1633.1666
$ws.Range("L73").Value = 4899.4998
$ws.Range("N73").Value = -6771.4998
$ws.Range("H129").Value = 1977.762
$ws.Range("J129").Value = 3971.5
$ws.Range("L129").Value = 11914.5
$ws.Range("N129").Value = -21914.5
$ws.Range("H132").Value = 1655.5349
$ws.Range("I132").Value = 1555.9744
$ws.Range("K132").Value = 4667.9232
$ws.Range("M132").Value = -2137.9232
$ws.Range("H136").Value = 78935.42999999999
$ws.Range("J136").Value = 78935.42999999999
$ws.Range("L136").Value = 78935.42999999999
$ws.Range("N136").Value = -89135.42999999999
$ws.Range("H138").Value = 2029.4286
$ws.Range("J138").Value = 2766.8845
$ws.Range("L138").Value = 8300.6535
$ws.Range("N138").Value = -18580.6535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3138.7183
$ws.Range("I32").Value = 2506.6924
$ws.Range("K32").Value = 2506.6924
$ws.Range("M32").Value = -2219.6924
$ws.Range("H45").Value = 11365718
$ws.Range("I45").Value = 2353.1428
$ws.Range("J45").Value = 31251606
$ws.Range("K45").Value = 2353.1428
$ws.Range("L45").Value = 31251606
$ws.Range("M45").Value = -1976.1428
$ws.Range("N45").Value = -31252360
$ws.Range("H61").Value = 51612.9
$ws.Range("I61").Value = 1265.3125
$ws.Range("K61").Value = 1265.3125
$ws.Range("M61").Value = -1053.3125
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H110").Value = 1214.5714
$ws.Range("I110").Value = 996.5
$ws.Range("K110").Value = 996.5
$ws.Range("M110").Value = 1048.5
$ws.Range("H122").Value = 2830.1738
$ws.Range("I122").Value = 2715
$ws.Range("K122").Value = 8145
$ws.Range("M122").Value = -5695
$ws.Range("H134").Value = 98408.5
$ws.Range("J134").Value = 98408.5
$ws.Range("L134").Value = 98408.5
$ws.Range("N134").Value = -108548.5
$ws.Range("H136").Value = 51612.9
$ws.Range("I136").Value = 1265.3125
$ws.Range("K136").Value = 3795.9375
$ws.Range("M136").Value = -1245.9375
$ws.Range("H139").Value = 149999
$ws.Range("J139").Value = 149999
$ws.Range("L139").Value = 149999
$ws.Range("N139").Value = -160279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1377649.8
$ws.Range("I7").Value = 3666833
$ws.Range("J7").Value = 4139.8
$ws.Range("K7").Value = 3666833
$ws.Range("L7").Value = 4139.8
$ws.Range("M7").Value = -3666720
$ws.Range("N7").Value = -4365.8
$ws.Range("H16").Value = 3581.6667
$ws.Range("I16").Value = 495
$ws.Range("J16").Value = 5125
$ws.Range("K16").Value = 495
$ws.Range("L16").Value = 5125
$ws.Range("M16").Value = -325
$ws.Range("N16").Value = -5465
$ws.Range("H94").Value = 1927.0416
$ws.Range("I94").Value = 1895.2222
$ws.Range("J94").Value = 2022.5
$ws.Range("K94").Value = 1895.2222
$ws.Range("L94").Value = 2022.5
$ws.Range("M94").Value = -1444.2222
$ws.Range("N94").Value = -2924.5
$ws.Range("H99").Value = 1113697.2
$ws.Range("I99").Value = 32939.03
$ws.Range("K99").Value = 32939.03
$ws.Range("M99").Value = -31441.03
$ws.Range("H105").Value = 41623.117
$ws.Range("I105").Value = 61173.59
$ws.Range("K105").Value = 61173.59
$ws.Range("M105").Value = -59426.59
$ws.Range("H134").Value = 3568.375
$ws.Range("I134").Value = 1232.6666
$ws.Range("J134").Value = 6571.4287
$ws.Range("K134").Value = 3697.9998
$ws.Range("L134").Value = 19714.2861
$ws.Range("M134").Value = -1162.9998
$ws.Range("N134").Value = -24784.2861
$ws.Range("H138").Value = 78815
$ws.Range("J138").Value = 79291.875
$ws.Range("L138").Value = 79291.875
$ws.Range("N138").Value = -89571.875
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 782
$ws.Range("I16").Value = 812.6667
$ws.Range("K16").Value = 812.6667
$ws.Range("M16").Value = -525.6667
$ws.Range("H58").Value = 1384.3448
$ws.Range("I58").Value = 1077.3636
$ws.Range("J58").Value = 2349.1428
$ws.Range("K58").Value = 1077.3636
$ws.Range("L58").Value = 2349.1428
$ws.Range("M58").Value = -874.3635999999999
$ws.Range("N58").Value = -2755.1428
$ws.Range("H74").Value = 64164.668
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 72997.60000000001
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 72997.60000000001
$ws.Range("M74").Value = -19126
$ws.Range("N74").Value = -74745.60000000001
$ws.Range("H77").Value = 64164.668
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 72997.60000000001
$ws.Range("K77").Value = 60000
$ws.Range("L77").Value = 218992.8
$ws.Range("M77").Value = -55632
$ws.Range("N77").Value = -227728.8
$ws.Range("H105").Value = 3248.4285
$ws.Range("I105").Value = 2257.8
$ws.Range("K105").Value = 2257.8
$ws.Range("M105").Value = -510.8000000000002
$ws.Range("H113").Value = 782
$ws.Range("I113").Value = 812.6667
$ws.Range("K113").Value = 812.6667
$ws.Range("M113").Value = 1357.3333
$ws.Range("H122").Value = 3983.818
$ws.Range("I122").Value = 3565.25
$ws.Range("K122").Value = 10695.75
$ws.Range("M122").Value = -8245.75
$ws.Range("H132").Value = 1959.2094
$ws.Range("I132").Value = 1802.2106
$ws.Range("J132").Value = 3152.4
$ws.Range("K132").Value = 5406.6318
$ws.Range("L132").Value = 9457.200000000001
$ws.Range("M132").Value = -2876.6318
$ws.Range("N132").Value = -14517.2
$ws.Range("H134").Value = 23191.797
$ws.Range("I134").Value = 2639.3171
$ws.Range("J134").Value = 128523.25
$ws.Range("K134").Value = 7917.951300000001
$ws.Range("L134").Value = 385569.75
$ws.Range("M134").Value = -5382.951300000001
$ws.Range("N134").Value = -390639.75
$ws.Range("H136").Value = 1384.3448
$ws.Range("I136").Value = 1077.3636
$ws.Range("J136").Value = 2349.1428
$ws.Range("K136").Value = 3232.0908
$ws.Range("L136").Value = 7047.428400000001
$ws.Range("M136").Value = -682.0907999999999
$ws.Range("N136").Value = -12147.4284
$ws.Range("H138").Value = 51706.668
$ws.Range("J138").Value = 49956.375
$ws.Range("L138").Value = 49956.375
$ws.Range("N138").Value = -60236.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2260
$ws.Range("I3").Value = 970
$ws.Range("K3").Value = 2910
$ws.Range("M3").Value = -2798
$ws.Range("H37").Value = 181192.31
$ws.Range("J37").Value = 181192.31
$ws.Range("L37").Value = 543576.9299999999
$ws.Range("N37").Value = -543800.9299999999
$ws.Range("H122").Value = 1011740.5
$ws.Range("I122").Value = 2449
$ws.Range("K122").Value = 22041
$ws.Range("M122").Value = -19591

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2047.68
$ws.Range("I97").Value = 1690
$ws.Range("K97").Value = 1690
$ws.Range("M97").Value = -1194
$ws.Range("H102").Value = 2691.2104
$ws.Range("I102").Value = 2720.6
$ws.Range("J102").Value = 2581
$ws.Range("K102").Value = 2720.6
$ws.Range("L102").Value = 2581
$ws.Range("M102").Value = -1098.6
$ws.Range("N102").Value = -5825
$ws.Range("H113").Value = 3335757.8
$ws.Range("I113").Value = 2437.5
$ws.Range("J113").Value = 5557971.5
$ws.Range("K113").Value = 2437.5
$ws.Range("L113").Value = 5557971.5
$ws.Range("M113").Value = -267.5
$ws.Range("N113").Value = -5562311.5
$ws.Range("H132").Value = 3080.8628
$ws.Range("I132").Value = 2583.025
$ws.Range("K132").Value = 7749.075000000001
$ws.Range("M132").Value = -5219.075000000001
$ws.Range("H138").Value = 139999.5
$ws.Range("J138").Value = 139999.5
$ws.Range("L138").Value = 139999.5
$ws.Range("N138").Value = -150279.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11145.381
$ws.Range("I7").Value = 12738.412
$ws.Range("J7").Value = 4375
$ws.Range("K7").Value = 12738.412
$ws.Range("L7").Value = 4375
$ws.Range("M7").Value = -12626.412
$ws.Range("N7").Value = -4599
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 4361.625
$ws.Range("I68").Value = 4882.1665
$ws.Range("J68").Value = 2800
$ws.Range("K68").Value = 4882.1665
$ws.Range("L68").Value = 2800
$ws.Range("M68").Value = -4133.1665
$ws.Range("N68").Value = -4298
$ws.Range("H71").Value = 4361.625
$ws.Range("I71").Value = 4882.1665
$ws.Range("J71").Value = 2800
$ws.Range("K71").Value = 24410.8325
$ws.Range("L71").Value = 14000
$ws.Range("M71").Value = -20666.8325
$ws.Range("N71").Value = -21488
$ws.Range("H82").Value = 3898.6667
$ws.Range("I82").Value = 4499
$ws.Range("J82").Value = 2998.1667
$ws.Range("K82").Value = 4499
$ws.Range("L82").Value = 2998.1667
$ws.Range("M82").Value = -4138
$ws.Range("N82").Value = -3720.1667
$ws.Range("H85").Value = 3898.6667
$ws.Range("I85").Value = 4499
$ws.Range("J85").Value = 2998.1667
$ws.Range("K85").Value = 4499
$ws.Range("L85").Value = 2998.1667
$ws.Range("M85").Value = -3251
$ws.Range("N85").Value = -5494.1667
$ws.Range("H93").Value = 1789.1
$ws.Range("I93").Value = 1862
$ws.Range("J93").Value = 1497.5
$ws.Range("K93").Value = 1862
$ws.Range("L93").Value = 1497.5
$ws.Range("M93").Value = -614
$ws.Range("N93").Value = -3993.5
$ws.Range("H100").Value = 11260.84
$ws.Range("I100").Value = 11919.611
$ws.Range("J100").Value = 9566.857
$ws.Range("K100").Value = 11919.611
$ws.Range("L100").Value = 9566.857
$ws.Range("M100").Value = -11378.611
$ws.Range("N100").Value = -10648.857
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 8713126
$ws.Range("J122").Value = 22227878
$ws.Range("L122").Value = 66683634
$ws.Range("N122").Value = -66688534
$ws.Range("H126").Value = 11145.381
$ws.Range("I126").Value = 12738.412
$ws.Range("J126").Value = 4375
$ws.Range("K126").Value = 38215.236
$ws.Range("L126").Value = 13125
$ws.Range("M126").Value = -35745.236
$ws.Range("N126").Value = -18065
$ws.Range("H132").Value = 1575
$ws.Range("I132").Value = 1279.5454
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 3838.6362
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -1308.6362
$ws.Range("N132").Value = -14660
$ws.Range("H134").Value = 139997.75
$ws.Range("J134").Value = 139997.75
$ws.Range("L134").Value = 139997.75
$ws.Range("N134").Value = -150137.75
$ws.Range("H138").Value = 143164.67
$ws.Range("J138").Value = 143164.67
$ws.Range("L138").Value = 143164.67
$ws.Range("N138").Value = -153444.67
$ws.Range("H141").Value = 93638
$ws.Range("J141").Value = 93638
$ws.Range("L141").Value = 93638
$ws.Range("N141").Value = -103998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7485
$ws.Range("I62").Value = 9681.777
$ws.Range("J62").Value = 5288.222
$ws.Range("K62").Value = 9681.777
$ws.Range("L62").Value = 5288.222
$ws.Range("M62").Value = -9057.777
$ws.Range("N62").Value = -6536.222
$ws.Range("H65").Value = 7485
$ws.Range("I65").Value = 9681.777
$ws.Range("J65").Value = 5288.222
$ws.Range("K65").Value = 48408.885
$ws.Range("L65").Value = 26441.11
$ws.Range("M65").Value = -45288.885
$ws.Range("N65").Value = -32681.11
$ws.Range("H96").Value = 2638793
$ws.Range("I96").Value = 8163.143
$ws.Range("K96").Value = 8163.143
$ws.Range("M96").Value = -6790.143
$ws.Range("H113").Value = 1304.8462
$ws.Range("I113").Value = 1400.4
$ws.Range("K113").Value = 4201.200000000001
$ws.Range("M113").Value = -2031.200000000001
$ws.Range("H126").Value = 3299.9443
$ws.Range("I126").Value = 3081.818
$ws.Range("J126").Value = 3642.7144
$ws.Range("K126").Value = 9245.454000000002
$ws.Range("L126").Value = 10928.1432
$ws.Range("M126").Value = -6775.454000000002
$ws.Range("N126").Value = -15868.1432
$ws.Range("H132").Value = 1554792.4
$ws.Range("I132").Value = 1919.4546
$ws.Range("J132").Value = 7248660
$ws.Range("K132").Value = 5758.3638
$ws.Range("L132").Value = 21745980
$ws.Range("M132").Value = -3228.3638
$ws.Range("N132").Value = -21751040
$ws.Range("H133").Value = 66698
$ws.Range("J133").Value = 64497.5
$ws.Range("L133").Value = 64497.5
$ws.Range("N133").Value = -74617.5
$ws.Range("H141").Value = 63500.74
$ws.Range("J141").Value = 63500.74
$ws.Range("L141").Value = 63500.74
$ws.Range("N141").Value = -73860.73999999999
